# Reorder the "Recorded By" (column G) entries per-row according to a fixed
# mapping of observed before/after values (the commit only reshuffles the
# comma-separated names/emails already present in each cell; no new data is
# introduced and no other columns are touched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of old cell text -> new cell text (built from the actual
# before/after values found in the workbook).
$map = @{
    "System, system, backup@backdoor.com" = "backup@backdoor.com, System, system"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

# Column G = "Recorded By" (7th column). Walk every data row (header is row 1).
$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 157) { $lastRow = 157 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($null -eq $val) { continue }
    $text = [string]$val
    if ($map.ContainsKey($text)) {
        $cell.Value = $map[$text]
    }
}
